# "append result of 刘思柒"
# Fill in the new student's row (row 16) in the results sheet, mirroring the
# pattern already used for the other rows (2-15): name/status pairs, a
# basic/advance/benchmark score triple, the measured time in J, and the
# rate/score/total formulas in K/L/M.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "515030910406_刘思柒"
$ws.Range("B16").Value = "Failed(WA)"
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = "Pass"
$ws.Range("G16").Value = 20
$ws.Range("H16").Value = "Failed(RE)"
$ws.Range("J16").Value = 2147483647

$ws.Range("K16").Formula = '=$J$17/J16'
$ws.Range("L16").Formula = '=20*K16'
$ws.Range("M16").Formula = '=D16+G16+L16'

# Matches the post-edit selection recorded in the sheet (user tabbed/entered
# past the last populated cell of the new row).
$ws.Range("H17").Select()
